$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 23, shifting rows 23+ down by one.
$ws.Rows.Item(23).Insert()

# Fill in the new row 23 data.
$ws.Range("A23").Value = 87082
$ws.Range("B23").Value = "https://leetcode.com/u/RayyanAshraf/"
$ws.Range("D23").Value = 287
$ws.Range("E23").Value = 15
$ws.Range("F23").Value = 58
$ws.Range("L23").Value = "https://github.com/etsryn"

# Update the style of A23 to match other rows in this block (style index 1 -> numFmtId 3, "#,##0").
$ws.Range("A23").NumberFormat = "#,##0"

# Update the active selection to L23 as per the diff.
$ws.Range("L23").Select()
